$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2200.8333
$ws.Range("I62").Value = 2241
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 2241
$ws.Range("L62").Value = 2000
$ws.Range("M62").Value = -1617
$ws.Range("N62").Value = -3248
$ws.Range("H65").Value = 2200.8333
$ws.Range("I65").Value = 2241
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 11205
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = -8085
$ws.Range("N65").Value = -16240
$ws.Range("H116").Value = 460015.97
$ws.Range("I116").Value = 716239.4
$ws.Range("J116").Value = 11625
$ws.Range("K116").Value = 716239.4
$ws.Range("L116").Value = 11625
$ws.Range("M116").Value = -712797.4
$ws.Range("N116").Value = -18509
$ws.Range("H132").Value = 368056.62
$ws.Range("I132").Value = 146515.7
$ws.Range("K132").Value = 439547.1
$ws.Range("M132").Value = -437017.1
$ws.Range("H137").Value = 563106.8
$ws.Range("I137").Value = 1703922.9
$ws.Range("J137").Value = 2705.9124
$ws.Range("K137").Value = 5111768.699999999
$ws.Range("L137").Value = 8117.7372
$ws.Range("M137").Value = -5109218.699999999
$ws.Range("N137").Value = -13217.7372
$ws.Range("H138").Value = 3326.5967
$ws.Range("I138").Value = 1860
$ws.Range("J138").Value = 3880.6445
$ws.Range("K138").Value = 5580
$ws.Range("L138").Value = 11641.9335
$ws.Range("M138").Value = -440
$ws.Range("N138").Value = -21921.9335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4382.4688
$ws.Range("I32").Value = 4664.7554
$ws.Range("J32").Value = 3713.8948
$ws.Range("K32").Value = 4664.7554
$ws.Range("L32").Value = 3713.8948
$ws.Range("M32").Value = -4377.7554
$ws.Range("N32").Value = -4287.8948
$ws.Range("H74").Value = 4315.5
$ws.Range("I74").Value = 4367
$ws.Range("K74").Value = 4367
$ws.Range("M74").Value = -3493
$ws.Range("H77").Value = 4315.5
$ws.Range("I77").Value = 4367
$ws.Range("K77").Value = 21835
$ws.Range("M77").Value = -17467
$ws.Range("H110").Value = 1430.6666
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 1430.6666
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 1430.6666
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = -5520.6666
$ws.Range("H132").Value = 1757.0238
$ws.Range("I132").Value = 984.19354
$ws.Range("J132").Value = 3935
$ws.Range("K132").Value = 2952.58062
$ws.Range("L132").Value = 11805
$ws.Range("M132").Value = -422.5806199999997
$ws.Range("N132").Value = -16865
$ws.Range("H137").Value = 44990
$ws.Range("J137").Value = 44990
$ws.Range("L137").Value = 44990
$ws.Range("N137").Value = -55190

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 34510
$ws.Range("J103").Value = 34510
$ws.Range("L103").Value = 34510
$ws.Range("N103").Value = -36854

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2315.1333
$ws.Range("I31").Value = 956.52
$ws.Range("J31").Value = 3285.5715
$ws.Range("K31").Value = 956.52
$ws.Range("L31").Value = 3285.5715
$ws.Range("M31").Value = -661.52
$ws.Range("N31").Value = -3875.5715
$ws.Range("H34").Value = 2315.1333
$ws.Range("I34").Value = 956.52
$ws.Range("J34").Value = 3285.5715
$ws.Range("K34").Value = 956.52
$ws.Range("L34").Value = 3285.5715
$ws.Range("M34").Value = -754.52
$ws.Range("N34").Value = -3689.5715
$ws.Range("H97").Value = 33746.125
$ws.Range("J97").Value = 33746.125
$ws.Range("L97").Value = 33746.125
$ws.Range("N97").Value = -35728.125
$ws.Range("H107").Value = 710.7
$ws.Range("I107").Value = 405.6207
$ws.Range("J107").Value = 1515
$ws.Range("K107").Value = 405.6207
$ws.Range("L107").Value = 1515
$ws.Range("M107").Value = 1514.3793
$ws.Range("N107").Value = -5355
$ws.Range("H138").Value = 48362.5
$ws.Range("J138").Value = 48362.5
$ws.Range("L138").Value = 48362.5
$ws.Range("N138").Value = -58642.5
$ws.Range("H140").Value = 80438.664
$ws.Range("J140").Value = 80438.664
$ws.Range("L140").Value = 80438.664
$ws.Range("N140").Value = -90798.664
$ws.Range("H141").Value = 29725
$ws.Range("J141").Value = 29725
$ws.Range("L141").Value = 29725
$ws.Range("N141").Value = -40085

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 418608.56
$ws.Range("I5").Value = 1020.4
$ws.Range("J5").Value = 787068.7
$ws.Range("K5").Value = 3061.2
$ws.Range("L5").Value = 2361206.1
$ws.Range("M5").Value = -2949.2
$ws.Range("N5").Value = -2361430.1
$ws.Range("H14").Value = 64.933334
$ws.Range("I14").Value = 64.933334
$ws.Range("K14").Value = 194.800002
$ws.Range("M14").Value = -21.80000200000001
$ws.Range("H68").Value = 1357.973
$ws.Range("I68").Value = 980.53845
$ws.Range("J68").Value = 1562.4166
$ws.Range("K68").Value = 2941.61535
$ws.Range("L68").Value = 4687.2498
$ws.Range("M68").Value = -2130.61535
$ws.Range("N68").Value = -6309.2498
$ws.Range("H71").Value = 1357.973
$ws.Range("I71").Value = 980.53845
$ws.Range("J71").Value = 1562.4166
$ws.Range("K71").Value = 8824.84605
$ws.Range("L71").Value = 14061.7494
$ws.Range("M71").Value = -4768.84605
$ws.Range("N71").Value = -22173.7494
$ws.Range("H80").Value = 27786174
$ws.Range("I80").Value = 7818.3335
$ws.Range("J80").Value = 33341846
$ws.Range("K80").Value = 23455.0005
$ws.Range("L80").Value = 100025538
$ws.Range("M80").Value = -22519.0005
$ws.Range("N80").Value = -100027410
$ws.Range("H83").Value = 27786174
$ws.Range("I83").Value = 7818.3335
$ws.Range("J83").Value = 33341846
$ws.Range("K83").Value = 70365.0015
$ws.Range("L83").Value = 300076614
$ws.Range("M83").Value = -65685.0015
$ws.Range("N83").Value = -300085974
$ws.Range("H129").Value = 1850.6
$ws.Range("J129").Value = 6650
$ws.Range("L129").Value = 19950
$ws.Range("N129").Value = -29950
$ws.Range("H131").Value = 913.83
$ws.Range("J131").Value = 991.1573
$ws.Range("L131").Value = 2973.4719
$ws.Range("N131").Value = -13053.4719
$ws.Range("H135").Value = 418608.56
$ws.Range("I135").Value = 1020.4
$ws.Range("J135").Value = 787068.7
$ws.Range("K135").Value = 9183.6
$ws.Range("L135").Value = 7083618.3
$ws.Range("M135").Value = -6648.6
$ws.Range("N135").Value = -7088688.3
$ws.Range("H137").Value = 3794.5625
$ws.Range("I137").Value = 4687.778
$ws.Range("J137").Value = 2646.1428
$ws.Range("K137").Value = 14063.334
$ws.Range("L137").Value = 7938.428400000001
$ws.Range("M137").Value = -8963.334000000001
$ws.Range("N137").Value = -18138.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3052.25
$ws.Range("I102").Value = 2085.4348
$ws.Range("J102").Value = 7499.6
$ws.Range("K102").Value = 2085.4348
$ws.Range("L102").Value = 7499.6
$ws.Range("M102").Value = -463.4348
$ws.Range("N102").Value = -10743.6
$ws.Range("H113").Value = 2203
$ws.Range("I113").Value = 2514.6667
$ws.Range("J113").Value = 800.5
$ws.Range("K113").Value = 2514.6667
$ws.Range("L113").Value = 800.5
$ws.Range("M113").Value = -344.6667000000002
$ws.Range("N113").Value = -5140.5
$ws.Range("H140").Value = 38930.555
$ws.Range("J140").Value = 38930.555
$ws.Range("L140").Value = 38930.555
$ws.Range("N140").Value = -49290.555

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1369.5294
$ws.Range("I61").Value = 1340.1666
$ws.Range("J61").Value = 1440
$ws.Range("K61").Value = 1340.1666
$ws.Range("L61").Value = 1440
$ws.Range("M61").Value = -1138.1666
$ws.Range("N61").Value = -1844
$ws.Range("H93").Value = 6947393.5
$ws.Range("I93").Value = 11113631
$ws.Range("J93").Value = 3663.3333
$ws.Range("K93").Value = 11113631
$ws.Range("L93").Value = 3663.3333
$ws.Range("M93").Value = -11112383
$ws.Range("N93").Value = -6159.3333
$ws.Range("H113").Value = 1369.5294
$ws.Range("I113").Value = 1340.1666
$ws.Range("J113").Value = 1440
$ws.Range("K113").Value = 1340.1666
$ws.Range("L113").Value = 1440
$ws.Range("M113").Value = 829.8334
$ws.Range("N113").Value = -5780
$ws.Range("H132").Value = 3149.1914
$ws.Range("I132").Value = 2394.8647
$ws.Range("J132").Value = 5940.2
$ws.Range("K132").Value = 7184.5941
$ws.Range("L132").Value = 17820.6
$ws.Range("M132").Value = -4654.5941
$ws.Range("N132").Value = -22880.6
$ws.Range("H136").Value = 3574.718
$ws.Range("I136").Value = 1583.75
$ws.Range("K136").Value = 4751.25
$ws.Range("M136").Value = -2201.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 224.57143
$ws.Range("I113").Value = 244.41176
$ws.Range("J113").Value = 193.90909
$ws.Range("K113").Value = 733.23528
$ws.Range("L113").Value = 581.72727
$ws.Range("M113").Value = 1436.76472
$ws.Range("N113").Value = -4921.72727
$ws.Range("H122").Value = 5310.4443
$ws.Range("J122").Value = 5936
$ws.Range("L122").Value = 17808
$ws.Range("N122").Value = -22708
$ws.Range("H132").Value = 10419034
$ws.Range("I132").Value = 1587.2
$ws.Range("K132").Value = 4761.6
$ws.Range("M132").Value = -2231.6
